# Add a new worksheet "FTNC_Average_Demand55" at the end of the workbook,
# mirroring the layout/content/style of the existing FTNC_Average_Demand5* sheets.

$wb = $excel.ActiveWorkbook

# Use the first existing sheet as the style/formatting template.
$template = $wb.Worksheets.Item(1)

# Insert the new sheet after the current last sheet, so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "FTNC_Average_Demand55"

# Header row
$ws.Range("B1").Value = "In-vehicle"
$ws.Range("C1").Value = "At-stop"
$ws.Range("D1").Value = "Extra"
$ws.Range("E1").Value = "Tardiness"
$ws.Range("F1").Value = "Total"

# Data row
$ws.Range("A2").Value = "FTNC_Average_Demand_5"
$ws.Range("B2").Value = 2323.045558379758
$ws.Range("C2").Value = 12927.3257339221
$ws.Range("D2").Value = 515.0616311181449
$ws.Range("E2").Value = 64.17206500977974
$ws.Range("F2").Value = 15829.60498842978

# Copy formatting (bold, centered, thin border) from the template sheet,
# reusing existing style entries rather than creating new ones.
$template.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
